# Weekly fruit/vegetable price update.
# A new daily price record (2022-07-04, serial 44746) is inserted above the
# existing row 350 on the single worksheet, pushing all subsequent rows
# (old 350..373) down by one (new 351..374). The workbook's used range
# grows from A1:R373 to A1:R374.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 350; Excel shifts rows 350:373 down to
# 351:374 and carries the existing row formatting (e.g. the date style on
# column D) onto the freshly inserted row.
$ws.Rows("350:350").Insert()

# Populate the newly inserted row 350 with the new price record.
$ws.Cells.Item(350, 1).Value  = 3
$ws.Cells.Item(350, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(350, 3).Value  = "Coquimbo"
$ws.Cells.Item(350, 4).Value  = 44746
$ws.Cells.Item(350, 5).Value  = 5
$ws.Cells.Item(350, 6).Value  = 100112043
$ws.Cells.Item(350, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(350, 8).Value  = "Sin especificar"
$ws.Cells.Item(350, 9).Value  = "Primera"
$ws.Cells.Item(350, 10).Value = 117
$ws.Cells.Item(350, 11).Value = 16000
$ws.Cells.Item(350, 12).Value = 17000
$ws.Cells.Item(350, 13).Value = 16513
$ws.Cells.Item(350, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(350, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(350, 16).Value = 236
$ws.Cells.Item(350, 17).Value = 70
$ws.Cells.Item(350, 18).Value = "Hortaliza"
